# Refresh the crypto price/volume table (Coin, Link, Price, Volume(1h))
# to match the latest scrape, per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '61.851.63'
$ws.Range('E2').Value = '  -0.90%  '

# Row 3
$ws.Range('D3').Value = '2.449.83'
$ws.Range('E3').Value = '  +0.83%  '

# Row 4
$ws.Range('E4').Value = '  +0.02%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '579.84'
$ws.Range('E5').Value = '  +0.89%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.27'
$ws.Range('E6').Value = '  -1.33%  '

# Row 7
$ws.Range('E7').Value = '  +0.14%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.528'
$ws.Range('E8').Value = '  -0.02%  '

# Row 9
$ws.Range('D9').Value = '2.444.30'
$ws.Range('E9').Value = '  +0.75%  '

# Row 10
$ws.Range('E10').Value = '  +1.95%  '

# Row 11
$ws.Range('E11').Value = '  +2.72%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.16'
$ws.Range('E12').Value = '  -1.07%  '

# Row 13
$ws.Range('E13').Value = '  -2.49%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.84'
$ws.Range('E14').Value = '  -1.90%  '

# Row 15
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '2.887.59'
$ws.Range('E15').Value = '  -0.13%  '

# Row 16
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000171'
$ws.Range('E16').Value = '  -1.03%  '

# Row 17
$ws.Range('D17').Value = '61.804.68'
$ws.Range('E17').Value = '  -1.01%  '

# Row 18
$ws.Range('D18').Value = '2.447.43'
$ws.Range('E18').Value = '  +0.70%  '

# Row 19
$ws.Range('E19').Value = '  -3.51%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.20'
$ws.Range('E20').Value = '  +1.43%  '

# Row 21
$ws.Range('B21').Value = 'LEO'
$ws.Range('C21').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.72'
$ws.Range('E21').Value = '  +12.34%  '

# Row 22
$ws.Range('B22').Value = 'BitcoinCash'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '324.79'
$ws.Range('E22').Value = '  -1.35%  '

# Row 23
$ws.Range('B23').Value = 'Polkadot'
$ws.Range('C23').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.07'
$ws.Range('E23').Value = '  -1.26%  '

# Row 24
$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.00'
$ws.Range('E24').Value = '  -0.10%  '

# Row 25
$ws.Range('B25').Value = 'SuiNetwork'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.92'
$ws.Range('E25').Value = '  -2.71%  '

# Row 26
$ws.Range('B26').Value = 'Litecoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '64.97'
$ws.Range('E26').Value = '  -1.16%  '

# Row 27
$ws.Range('B27').Value = 'Aptos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.11'
$ws.Range('E27').Value = '  +2.47%  '

# Row 28
$ws.Range('B28').Value = 'Bittensor'
$ws.Range('C28').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '584.28'
$ws.Range('E28').Value = '  -7.51%  '

# Row 29
$ws.Range('B29').Value = 'WrappedeETH'
$ws.Range('C29').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D29').Value = '2.562.69'
$ws.Range('E29').Value = '  +0.02%  '

# Row 30
$ws.Range('B30').Value = 'Binance-PegBSC-USD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  -0.13%  '

# Row 31
$ws.Range('B31').Value = 'PEPE'
$ws.Range('C31').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D31').Value = '0.0₃0925'
$ws.Range('E31').Value = '  -3.48%  '

# Row 32
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.87'
$ws.Range('E32').Value = '  -2.00%  '

# Row 33
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.37'
$ws.Range('E33').Value = '  -4.33%  '

# Row 34
$ws.Range('B34').Value = 'PancakeSwap'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.87'
$ws.Range('E34').Value = '  -0.09%  '

# Row 35
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.134'
$ws.Range('E35').Value = '  -3.17%  '

# Row 36
$ws.Range('B36').Value = 'FirstDigitalUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  +0.04%  '

# Row 37
$ws.Range('B37').Value = 'NEARProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.75'
$ws.Range('E37').Value = '  -4.38%  '

# Row 38
$ws.Range('B38').Value = 'PolygonEcosystemToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.371'
$ws.Range('E38').Value = '  -1.06%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '151.84'
$ws.Range('E39').Value = '  +2.17%  '

# Row 40
$ws.Range('B40').Value = 'ImmutableX'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.40'
$ws.Range('E40').Value = '  -4.14%  '

# Row 41
$ws.Range('B41').Value = 'EthereumClassic'
$ws.Range('C41').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '18.28'
$ws.Range('E41').Value = '  -0.94%  '

# Row 42
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.14'
$ws.Range('E42').Value = '  -2.43%  '

# Row 43
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '42.22'
$ws.Range('E43').Value = '  -0.33%  '

# Row 44
$ws.Range('B44').Value = 'USDe'
$ws.Range('C44').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.999'
$ws.Range('E44').Value = '  +0.01%  '

# Row 45
$ws.Range('B45').Value = 'Stacks'
$ws.Range('C45').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.67'
$ws.Range('E45').Value = '  -3.90%  '

# Row 46
$ws.Range('B46').Value = 'dogwifhat'
$ws.Range('C46').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.34'
$ws.Range('E46').Value = '  -5.05%  '

# Row 47
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').Value = '0.0₆0276'
$ws.Range('E47').Value = '  +18.93%  '

# Row 48
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '141.15'
$ws.Range('E48').Value = '  -2.23%  '

# Row 49
$ws.Range('B49').Value = 'Filecoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.56'
$ws.Range('E49').Value = '  -3.42%  '

# Row 50
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.597'
$ws.Range('E50').Value = '  +0.38%  '

# Row 51
$ws.Range('B51').Value = 'Hedera'
$ws.Range('C51').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0509'
$ws.Range('E51').Value = '  -2.31%  '
